$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "Categoria" - copy formatting (bold/border/centered)
# from the existing header cell A1, then set its text.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Categoria"

# Update the "Período" column (C) values and populate the new "Categoria"
# column (D) for every data row.
$ws.Range("C2").Value = "2010 / 2022"
$ws.Range("D2").Value = "Toda a série histórica"

$ws.Range("C3").Value = "2010 / 2022"
$ws.Range("D3").Value = "Toda a série histórica"

$ws.Range("C4").Value = "2010 / 2022"
$ws.Range("D4").Value = "Toda a série histórica"

$ws.Range("C5").Value = "2010 / 2022"
$ws.Range("D5").Value = "Toda a série histórica"

$ws.Range("C6").Value = "2010 / 2022"
$ws.Range("D6").Value = "Toda a série histórica"

$ws.Range("C7").Value = "2010 / 2022"
$ws.Range("D7").Value = "Toda a série histórica"

$ws.Range("C8").Value = "2021 / 2022"
$ws.Range("D8").Value = "Último ano da série histórica"

$ws.Range("C9").Value = "2021 / 2022"
$ws.Range("D9").Value = "Último ano da série histórica"

$ws.Range("C10").Value = "2021 / 2022"
$ws.Range("D10").Value = "Último ano da série histórica"

$ws.Range("C11").Value = "2021 / 2022"
$ws.Range("D11").Value = "Último ano da série histórica"

$ws.Range("C12").Value = "2021 / 2022"
$ws.Range("D12").Value = "Último ano da série histórica"

$ws.Range("C13").Value = "2021 / 2022"
$ws.Range("D13").Value = "Último ano da série histórica"
